$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71:90 down to 72:91.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 45120
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = 100112022
$ws.Range("G71").Value = "Arveja Verde"
$ws.Range("H71").Value = "Perfection"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 35
$ws.Range("K71").Value = 23000
$ws.Range("L71").Value = 23000
$ws.Range("M71").Value = 23000
$ws.Range("N71").Value = "$/saco 25 kilos"
$ws.Range("O71").Value = "Provincia de Limarí"
$ws.Range("P71").Value = 920
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"
